# Add the new "Git link of the assignment an assessment" column header
# next to the existing TRAINER RATING header in row 6.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = "Git link of the assignment an assessment"

# Widen column E so the new header text fits (previously it was a
# narrow, bestFit column sized for something else).
$ws.Columns.Item(5).ColumnWidth = 47.72

# Keep the active selection on the newly added cell, mirroring the
# original author's last selection when they saved the workbook.
$ws.Range("E6").Select()
